$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USER8")

$ws.Columns.Item(3).ColumnWidth = 2.88

$ws.Range("B2").Value = 0.27826086956521739
$ws.Range("C2").Value = 0.20000000000000001
$ws.Range("D2").Value = 0.63478260869565217
$ws.Range("E2").Value = 0.66666666666666663

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0.16666666666666666
$ws.Range("C3").Value = 0.10000000000000001
$ws.Range("D3").Value = 0.16666666666666666
$ws.Range("E3").Value = 0.5
